$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# NOTE: the fill order below is deliberate (not alphabetical / not
# left-to-right) so that new shared-string table entries get appended in
# the same order Excel produced them: mailTemplate, mailTableName,
# mailRowName, templateOfMandrill, default, MyTemplate.

# Row 18: mailTableName
$ws.Range("A18").Value = "*"
$ws.Range("B18").Value = "mandrill"
$ws.Range("E18").Value = "mailTemplate"
$ws.Range("F18").Value = "mailTemplate"
$ws.Range("G18").Value = "mailTemplate"
$ws.Range("C18").Value = "mailTableName"
$ws.Range("D18").Formula = '="@@."&A18&"."&B18&"."&C18&"@@"'

# Row 19: mailRowName
$ws.Range("A19").Value = "*"
$ws.Range("B19").Value = "mandrill"
$ws.Range("C19").Value = "mailRowName"
$ws.Range("D19").Formula = '="@@."&A19&"."&B19&"."&C19&"@@"'

# Row 20: templateOfMandrill
$ws.Range("A20").Value = "*"
$ws.Range("B20").Value = "mandrill"
$ws.Range("C20").Value = "templateOfMandrill"
$ws.Range("D20").Formula = '="@@."&A20&"."&B20&"."&C20&"@@"'

$ws.Range("E19").Value = "default"
$ws.Range("F19").Value = "default"
$ws.Range("G19").Value = "default"

$ws.Range("E20").Value = "MyTemplate"
$ws.Range("F20").Value = "MyTemplate"
$ws.Range("G20").Value = "MyTemplate"

# Update view state to match the new selection/scroll position
$ws.Range("F19").Select()
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 2
